$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (FE)
$ws.Range("E2").Value = 0.44
$ws.Range("F2").Value = 0.9
$ws.Range("G2").Value = 0.7

# Row 3 (FE+Disg) - fill previously missing values
$ws.Range("E3").Value = 0.37
$ws.Range("F3").Value = 0.9
$ws.Range("G3").Value = 0.7

# Row 4 (FE+Disg+Var)
$ws.Range("B4").Value = 0.34
$ws.Range("E4").Value = 0.55
$ws.Range("F4").Value = 0.9
$ws.Range("G4").Value = 0.8100000000000001
